$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, shifting existing rows 80:202 down to 81:203
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with the new weekly price record
$ws.Range("A80").Value = 11
$ws.Range("B80").Value = "Vega Monumental Concepción"
$ws.Range("C80").Value = "Bíobío"
$ws.Range("D80").Value = 45036
$ws.Range("E80").Value = 8
$ws.Range("F80").Value = 100112043
$ws.Range("G80").Value = "Pepino ensalada"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 100
$ws.Range("K80").Value = 10000
$ws.Range("L80").Value = 11000
$ws.Range("M80").Value = 10500
$ws.Range("N80").Value = "`$/caja 60 unidades"
$ws.Range("O80").Value = "Región de Arica y Parinacota"
$ws.Range("P80").Value = 175
$ws.Range("Q80").Value = 60
$ws.Range("R80").Value = "Hortaliza"
